$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the cells below can be
# updated, then restore protection at the end.
$ws.Unprotect()

# Bump the "as of" date in the confidential disclosure footnote from
# 2021-04-05 to 2021-04-06.
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-06 for illustrative purposes only and are subject to change."

# Setting a multi-line value can make Excel mark the row with an explicit
# custom height; AutoFit puts it back to the sheet's default row height
# behavior (no custom height flag), matching the original row.
$ws.Rows.Item(9).AutoFit()

# Refresh the Weight (D) and Percent Change (E) columns for rows 2-6.
$ws.Range("D2").Value = 0.2503270204826472
$ws.Range("E2").Value = 0.002236858456567736

$ws.Range("D3").Value = 0.2464062227527315
$ws.Range("E3").Value = -0.002014968336211953

$ws.Range("D4").Value = 0.2550558773967936
$ws.Range("E4").Value = -0.004266705235753565

$ws.Range("D5").Value = 0.2482108793678276
$ws.Range("E5").Value = 0.002495075508863964

$ws.Range("D6").Value = 0.9999999999999999
$ws.Range("E6").Value = -0.0004054979853728158

# Restore sheet protection so the workbook stays locked for end users.
$ws.Protect()
